$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (L..O). Insert shared strings in the same order the
# original author's diff lists them: GC, Genes, Coding, Size bp.
$ws.Range("M1").Value = "GC"
$ws.Range("N1").Value = "Genes"
$ws.Range("O1").Value = "Coding"
$ws.Range("L1").Value = "Size bp"

# Row 2 - FDAARGOS_1054 / ATCC 6940
$ws.Range("L2").Value = 2952500
$ws.Range("M2").Value = 59
$ws.Range("N2").Value = 2783
$ws.Range("O2").Value = 2608

# Row 3 - FDAARGOS_1197
$ws.Range("L3").Value = 2948781
$ws.Range("M3").Value = 59
$ws.Range("N3").Value = 2830
$ws.Range("O3").Value = 2830

# Row 4 - FDAARGOS_1115
$ws.Range("L4").Value = 2904831
$ws.Range("M4").Value = 59
$ws.Range("N4").Value = 2716
$ws.Range("O4").Value = 2569

# Row 5 - FDAARGOS_1116
$ws.Range("L5").Value = 2665682
$ws.Range("M5").Value = 59
$ws.Range("N5").Value = 2471
$ws.Range("O5").Value = 2345

# Row 7 - KC-Na-01
$ws.Range("L7").Value = 2758500
$ws.Range("M7").Value = 59
$ws.Range("N7").Value = 2653
$ws.Range("O7").Value = 2517

# Apply thousands-separator number format to the genome-size / gene-count /
# coding-count columns (not the GC% column).
$ws.Range("L2:L5").NumberFormat = "#,##0"
$ws.Range("L7").NumberFormat = "#,##0"
$ws.Range("N2:N5").NumberFormat = "#,##0"
$ws.Range("N7").NumberFormat = "#,##0"
$ws.Range("O2:O5").NumberFormat = "#,##0"
$ws.Range("O7").NumberFormat = "#,##0"

# Move the active selection, matching the author's final cursor position.
$ws.Range("C29").Select() | Out-Null
